$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 11; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "Videnov"
    $ws.Cells.Item($r, 2).Value = "Sofia, Tsarigradsko, 15"
}
